# Fix a numbering gap in column E (the "No." sequence within each group
# defined by column D). Several groups were missing one sequence number;
# every cell below the gap shifts up by one position, closing the gap.
# Column E holds text (shared-string) values such as "3", "4", "5", ... -
# the leading apostrophe forces each assignment to be stored as text
# (matching the existing text values in that column) instead of being
# auto-converted to a number.
#
# Block 1: group D=4,  rows 29-36  -> E sequence was missing "4"; shift
#          rows 32-36 up by one slot (5,6,7,8,9 -> 4,5,6,7,8).
# Block 2: group D=9,  rows 75-83  -> E sequence was missing "5"; shift
#          rows 79-83 up by one slot (6,7,8,9,10 -> 5,6,7,8,9).
# Block 3: group D=22, rows 149-154 -> E sequence was missing "3"; shift
#          rows 151-154 up by one slot (4,5,6,7 -> 3,4,5,6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E32").Value = "'4"
$ws.Range("E33").Value = "'5"
$ws.Range("E34").Value = "'6"
$ws.Range("E35").Value = "'7"
$ws.Range("E36").Value = "'8"

$ws.Range("E79").Value = "'5"
$ws.Range("E80").Value = "'6"
$ws.Range("E81").Value = "'7"
$ws.Range("E82").Value = "'8"
$ws.Range("E83").Value = "'9"

$ws.Range("E151").Value = "'3"
$ws.Range("E152").Value = "'4"
$ws.Range("E153").Value = "'5"
$ws.Range("E154").Value = "'6"
